# Auto update Excel log 2026-02-04 14:24:46
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets.

$wb = $excel.ActiveWorkbook

$pirRows = @(
    ,@(270, "14:23:42", "14:00", "Bathroom", "Motion Detected", "Active")
    ,@(271, "14:23:44", "14:00", "Bathroom", "No Motion", "Inactive")
    ,@(272, "14:23:46", "14:00", "Bathroom", "Motion Detected", "Active")
    ,@(273, "14:23:52", "14:00", "Bathroom", "No Motion", "Inactive")
    ,@(274, "14:23:58", "14:00", "Bathroom", "No Motion", "Inactive")
    ,@(275, "14:24:03", "14:00", "Bathroom", "No Motion", "Inactive")
    ,@(276, "14:24:08", "14:00", "Bathroom", "No Motion", "Inactive")
    ,@(277, "14:24:13", "14:00", "Bathroom", "No Motion", "Inactive")
    ,@(278, "14:24:18", "14:00", "Bathroom", "No Motion", "Inactive")
    ,@(279, "14:24:21", "14:00", "Bathroom", "Motion Detected", "Active")
    ,@(280, "14:24:27", "14:00", "Bathroom", "No Motion", "Inactive")
    ,@(281, "14:24:28", "14:00", "Bathroom", "Motion Detected", "Active")
    ,@(282, "14:24:36", "14:00", "Bathroom", "No Motion", "Inactive")
    ,@(283, "14:24:38", "14:00", "Bathroom", "Motion Detected", "Active")
)

$humidityRows = @(
    ,@(224, "14:23:40", "14:00", "Bathroom", "78.8%", "Active")
    ,@(225, "14:23:42", "14:00", "Bathroom", "78.0%", "Active")
    ,@(226, "14:23:45", "14:00", "Bathroom", "79.0%", "Active")
    ,@(227, "14:23:50", "14:00", "Bathroom", "78.2%", "Active")
    ,@(228, "14:23:55", "14:00", "Bathroom", "79.2%", "Active")
    ,@(229, "14:24:00", "14:00", "Bathroom", "78.5%", "Active")
    ,@(230, "14:24:05", "14:00", "Bathroom", "79.4%", "Active")
    ,@(231, "14:24:10", "14:00", "Bathroom", "78.5%", "Active")
    ,@(232, "14:24:15", "14:00", "Bathroom", "79.5%", "Active")
    ,@(233, "14:24:20", "14:00", "Bathroom", "78.5%", "Active")
    ,@(234, "14:24:25", "14:00", "Bathroom", "79.6%", "Active")
    ,@(235, "14:24:30", "14:00", "Bathroom", "78.7%", "Active")
    ,@(236, "14:24:40", "14:00", "Bathroom", "78.7%", "Active")
)

$temperatureRows = @(
    ,@(224, "14:23:41", "14:00", "Bathroom", "24.5C", "Active")
    ,@(225, "14:23:43", "14:00", "Bathroom", "24.5C", "Active")
    ,@(226, "14:23:45", "14:00", "Bathroom", "24.5C", "Active")
    ,@(227, "14:23:50", "14:00", "Bathroom", "24.4C", "Active")
    ,@(228, "14:23:55", "14:00", "Bathroom", "24.4C", "Active")
    ,@(229, "14:24:00", "14:00", "Bathroom", "24.4C", "Active")
    ,@(230, "14:24:05", "14:00", "Bathroom", "24.3C", "Active")
    ,@(231, "14:24:10", "14:00", "Bathroom", "24.3C", "Active")
    ,@(232, "14:24:15", "14:00", "Bathroom", "24.3C", "Active")
    ,@(233, "14:24:20", "14:00", "Bathroom", "24.3C", "Active")
    ,@(234, "14:24:26", "14:00", "Bathroom", "24.3C", "Active")
    ,@(235, "14:24:31", "14:00", "Bathroom", "24.3C", "Active")
)

# ---------------------------------------------------------------------
# PIR sheet: columns A Date, B Timestamp, C Hour, D Location, E Value, F Status
# ---------------------------------------------------------------------
$wsPir = $wb.Worksheets.Item("PIR")
$wsPir.Range("A270:A283").NumberFormat = "@"
foreach ($r in $pirRows) {
    $rowNum = $r[0]
    $wsPir.Cells.Item($rowNum, 1).Value = "2026-02-04"
    $wsPir.Cells.Item($rowNum, 2).Value = $r[1]
    $wsPir.Cells.Item($rowNum, 3).Value = $r[2]
    $wsPir.Cells.Item($rowNum, 4).Value = $r[3]
    $wsPir.Cells.Item($rowNum, 5).Value = $r[4]
    $wsPir.Cells.Item($rowNum, 6).Value = $r[5]
}

# ---------------------------------------------------------------------
# Humidity sheet: Value column (E) holds percentages, e.g. "78.8%"
# ---------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")
$wsHumidity.Range("A224:A236").NumberFormat = "@"
$wsHumidity.Range("E224:E236").NumberFormat = "@"
foreach ($r in $humidityRows) {
    $rowNum = $r[0]
    $wsHumidity.Cells.Item($rowNum, 1).Value = "2026-02-04"
    $wsHumidity.Cells.Item($rowNum, 2).Value = $r[1]
    $wsHumidity.Cells.Item($rowNum, 3).Value = $r[2]
    $wsHumidity.Cells.Item($rowNum, 4).Value = $r[3]
    $wsHumidity.Cells.Item($rowNum, 5).Value = $r[4]
    $wsHumidity.Cells.Item($rowNum, 6).Value = $r[5]
}

# ---------------------------------------------------------------------
# Temperature sheet: Value column (E) holds values like "24.5C"
# ---------------------------------------------------------------------
$wsTemperature = $wb.Worksheets.Item("Temperature")
$wsTemperature.Range("A224:A235").NumberFormat = "@"
foreach ($r in $temperatureRows) {
    $rowNum = $r[0]
    $wsTemperature.Cells.Item($rowNum, 1).Value = "2026-02-04"
    $wsTemperature.Cells.Item($rowNum, 2).Value = $r[1]
    $wsTemperature.Cells.Item($rowNum, 3).Value = $r[2]
    $wsTemperature.Cells.Item($rowNum, 4).Value = $r[3]
    $wsTemperature.Cells.Item($rowNum, 5).Value = $r[4]
    $wsTemperature.Cells.Item($rowNum, 6).Value = $r[5]
}

Write-Host "PIR rows now:" $wsPir.UsedRange.Rows.Count
Write-Host "Humidity rows now:" $wsHumidity.UsedRange.Rows.Count
Write-Host "Temperature rows now:" $wsTemperature.UsedRange.Rows.Count
